$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B header text changes from cDNA_ug to cDNA_ng
$ws.Range("B1").Value = "cDNA_ng"

# Column B values change from 1 to 1000 for rows 2-17
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 2).Value = 1000
}

# Update the selection to match the recorded view state
$ws.Range("F15").Select()
